$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that received updated values (row-wise), in order: B,D,E,F,G,I,K,L,M
$cols = @(2,4,5,6,7,9,11,12,13)   # column indices for B,D,E,F,G,I,K,L,M

$data = @(
    @(12.48646723298373,8.498892129854536,11.35087430413469,43.74201707914073,3.724520570734143,34.82025943552465,12.91836874913834,10.70111164754587,14.96453792705237),
    @(12.47837265242513,8.490878560581223,11.26040879204615,43.31310791699558,3.728176500998455,34.70402485822505,12.71931627174165,10.69179746896463,14.97896920775743),
    @(12.4767950400259,8.486089419064676,11.2032005886357,43.05554038489529,3.730537402822632,34.63565688092577,12.60046836262643,10.68813308339531,14.99062752729272),
    @(12.47700819603486,8.484170644024026,11.17947241101096,42.95212044301346,3.731528805196999,34.60855937001526,12.55295551044128,10.68715763106998,14.99608226175367),
    @(12.47709536793845,8.483854011391736,11.17550735329099,42.93504316551072,3.731695200630179,34.60410616508326,12.54512355595815,10.68702696108243,14.99703053734345),
    @(12.47679444518239,8.486063409060812,11.20288226046596,43.05413927687277,3.730550654379396,34.63528833454121,12.59982377589453,10.68811783018272,14.99069824138167),
    @(12.48297408507187,8.496101722575439,11.32002825451775,43.59298144351192,3.725757095699801,34.77956429319817,12.84907887990602,10.69747449097552,14.96893353427108),
    @(12.52185148840035,8.516835562254293,11.53645932387504,44.69161015597999,3.717273493444726,35.08589912449098,13.3614243273703,10.73206125181501,14.948424238932),
    @(12.56647419709497,8.532719115600734,11.68723998664041,45.51892339074104,3.711592327558877,35.32466733871539,13.74788698701153,10.76726896992274,14.9468286829624),
    @(12.59019391443621,8.540088349329835,11.75402102728206,45.89847777145178,3.709126123648288,35.43612993521457,13.92497049498811,10.7853848664779,14.94901476446358),
    @(12.59966112478821,8.542899668773314,11.77904671846042,46.04256516992868,3.70820911650369,35.47873437298021,13.99213748959087,10.7925438710178,14.95025975413245),
    @(12.59760073316706,8.542293276777194,11.77366871834788,46.01151906183748,3.708405860836038,35.46954135157557,13.97766817994842,10.79098880887826,14.94997309693928),
    @(12.59096310726568,8.540319221978059,11.75608518819109,45.91032533253394,3.709050342999568,35.43962717018206,13.93049459342219,10.78596787000685,14.94910883959332),
    @(12.58696032715416,8.539112760957046,11.74528046333663,45.84838488354794,3.709447303517305,35.42135503371733,13.90161146503229,10.78293123128125,14.94863373440216),
    @(12.56499232127252,8.532240418407602,11.68283904712533,45.49417435370333,3.711755869074175,35.31743914198552,13.73633348039998,10.76612705348606,14.94674428420085),
    @(12.55238711573078,8.528061334215408,11.64406785321729,45.27762569795395,3.713202295358495,35.2544097518729,13.63521475487985,10.75635396922506,14.94632993880716),
    @(12.54545934816838,8.525671159654976,11.62159748740703,45.15338050234989,3.714045372030598,35.21842560223515,13.57717874577735,10.75093058765435,14.94636580869688),
    @(12.54316929723723,8.524864207623462,11.61396025563998,45.11136921465301,3.714332738171685,35.20628852792245,13.55755249608528,10.7491283898809,14.94642508989233),
    @(12.55369564026757,8.528504804532727,11.64821274692762,45.30064655269512,3.713047169605842,35.2610916238882,13.64596664670458,10.75737387624973,14.94634567956174),
    @(12.59289963164439,8.54089848597156,11.76125705085735,45.94003947356855,3.708860585384815,35.44840305247934,13.94434823712719,10.78743455489709,14.94935138467),
    @(12.62134546257904,8.54911945143807,11.83360486413971,46.35996135902123,3.706222814685615,35.57312334634673,14.13996200061882,10.80882158390842,14.9537464766175),
    @(12.60590742294234,8.544720674356258,11.79513259453649,46.13568814375699,3.707621672939561,35.50635173095547,14.03552807016768,10.79724873720102,14.95117889943092),
    @(12.5531030617594,8.528304272862599,11.64633940192948,45.29023803293407,3.713117266107514,35.25806996264838,13.64110540389705,10.75691216773557,14.94633770931505),
    @(12.50849296552688,8.511113033580701,11.47934539709396,44.39046295600446,3.719471136668233,35.00057690085976,13.22074559067002,10.72097550099727,14.95160260992187)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($rowNum, $cols[$j]).Value = $rowVals[$j]
    }
}

$wb.Save()

